# Implemented track generation from GPS data for bahrein
#
# Updates car-setup parameters that were re-derived for the Bahrein track
# (from GPS-based track generation) and leaves each touched sheet's
# selection where the user's cursor ended up while reviewing the change.

$wb = $excel.ActiveWorkbook

# --- Engine: lower final drive / max power for the Bahrein setup ---------
$wsEngine = $wb.Worksheets.Item("Engine")
$wsEngine.Activate() | Out-Null
$wsEngine.Range("A2").Value = 2.9     # final_gear_ratio: 3.1 -> 2.9
$wsEngine.Range("C2").Value = 500     # maximum_power: 550 -> 500
$wsEngine.Range("C3").Select() | Out-Null

# --- Aero: reduced rear downforce coefficient -----------------------------
$wsAero = $wb.Worksheets.Item("Aero")
$wsAero.Activate() | Out-Null
$wsAero.Range("B2").Value = 4.5       # C_down: 5 -> 4.5
$wsAero.Range("B3").Select() | Out-Null

# --- Susp: raised lateral load transfer distribution ----------------------
$wsSusp = $wb.Worksheets.Item("Susp")
$wsSusp.Activate() | Out-Null
$wsSusp.Range("A2").Value = 0.55      # LLTD: 0.51 -> 0.55

# --- TireRear: just repositioning the cursor while reviewing --------------
$wsTireRear = $wb.Worksheets.Item("TireRear")
$wsTireRear.Activate() | Out-Null
$wsTireRear.Range("G25").Select() | Out-Null

# --- MassInertia: heavier sprung mass / yaw inertia for new spec ----------
$wsMass = $wb.Worksheets.Item("MassInertia")
$wsMass.Activate() | Out-Null
$wsMass.Range("C2").Value = 1050      # Ms: 965 -> 1050
$wsMass.Range("D2").Value = 1750      # Iz: 1400 -> 1750
$wsMass.Range("D3").Select() | Out-Null
